$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 272
$ws1.Range("F5").Value = 153
$ws1.Range("F6").Value = 55
$ws1.Range("F7").Value = 268
$ws1.Range("F8").Value = 214
$ws1.Range("F9").Value = 1985
$ws1.Range("F11").Value = 4672
$ws1.Range("F12").Value = 84
$ws1.Range("F13").Value = 329

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 272
$ws4.Range("F7").Value = 153
$ws4.Range("F8").Value = 55
$ws4.Range("F9").Value = 268
$ws4.Range("F10").Value = 214
$ws4.Range("F13").Value = 1985
$ws4.Range("F15").Value = 4672
$ws4.Range("F16").Value = 84
$ws4.Range("F17").Value = 329
